$d = $word.ActiveDocument

function AddParaBeforeEnd($text) {
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)
    $insertPos = $lastPara.Range.Start
    $r = $d.Range($insertPos, $insertPos)
    $r.InsertBefore($text + "`r")
}

AddParaBeforeEnd("27/10 Install OpenCard framework and read documentation")
AddParaBeforeEnd("28/10 Started working on skeleton code for host-side app using code samples and OpenCard api docs.")
AddParaBeforeEnd("29/10 Wrote most of both host-side and card-side app, haven’t yet worked out how to install onto card because I don’t have my card reader.")
AddParaBeforeEnd("31/10 Reader arrived. Read more of the documentation, built a project. Installed gpshell, installed reader driver, successfully loaded gpshell sample applet onto card, listed contents, and removed it.")
AddParaBeforeEnd("2/11 Attempted to load project onto card, discovered that my java cards are not JCDK 3.0.5 compatible. Choice: Get more up-to-date java cards, or use the older build process with the older JCDK? Reverting to older process means no handy Eclipse plugin. Documentation is also much better in newer version. Decided to buy new cards.")
AddParaBeforeEnd("4/11 Attempted to run sample project by running OCF host app and connecting to card. OCF wouldn’t recognise the card. I think I’ve identified that the reader needs an OCF driver, but it only has a PC/SC driver. Sources are few and far between. May have to redo application in C++ to use PCSC API. Could also be incompatibility with Windows 10, perhaps try in Ubuntu first.")
AddParaBeforeEnd("Data sheet for reader only specifies PCSC whereas for some other readers their data sheets also mention OCF. ")
AddParaBeforeEnd("From OCF programmer’s guide: “The OpenCard Reference implementation comes with a lockable CardTerminal implementation for PCSC card readers”. Downloaded reference implementation, there’s source code for PCSC-related stuff. Pcsc-wrapper-src.jar contains class files. Should look into usage. Very little information online about them. In particular, contains Pcsc10CardTerminal which emulates OCF stuff on PCSC.")
AddParaBeforeEnd("New cards arrived. Tried installing apps, got Unknown ISO7816 error: 0x6438 for apps compiled via eclipse. Possibly because card is jcdk 3.0.4 compatible, compiled using jcdk 3.0.5. Appears to be card-defined execution error.")
AddParaBeforeEnd("5/11 Decided to give the Python library pyscard a try instead of OCF. Spent the day running into various problems to do with dependencies it couldn’t find. Turns out plugin only for python 2, but I was trying python3. Working now, but unable to establish context. Not sure what the problem was, but it appears to be that particular version (1.7.0), because it worked when I tried a different version. With it, I was able to successfully test a small program that sends a SELECT APDU to the ISD and prints the sw1-sw2 output. ")
AddParaBeforeEnd("The pyscard Python extension is likely to be my tool of choice going forward. It’s compatible with both Microsoft PC/SC and Linux PC/SC lite, so the code will be portable and relatively easy to write.")
AddParaBeforeEnd("7/11 Followed through applet compiling process for JCDK 2.2.2. Problem running first demo with jcwde, gives message “card was unexpected at this time”. No idea what it means. TODO: Try another demo. Not a serious problem though, it’s the simulation test suite so it’s optional.")
AddParaBeforeEnd("Compiled sample code to class file, attempted to convert to CAP file. Converter has message “card was unexpected at this time”. Wtf?")
AddParaBeforeEnd("Identified. It doesn’t like spaces in the JC_HOME environment variable. Had to change directory structure.")
AddParaBeforeEnd("New problem: Script also can’t deal with space in path of JAVA_HOME. Scripts involve something like %JAVA_HOME%\bin\java -classpath %_CLASSES% com.sun.javacard.converter.Converter %*, where they should have quotes around the environment variable. Error on distributor’s part. Have to reinstall jdk into a different directory tomorrow.")
AddParaBeforeEnd("08/11 Was able to solve the problems and compile and convert an applet, and store/remove it from the card. Used JCDK HelloWorld sample applet source file.")
AddParaBeforeEnd("09/11 Looked into different asymmetric cryptography protocols.")
AddParaBeforeEnd("11/11 Configured Atom with Python and Java IDEs so I could develop both ends side-by-side. Adapted the build process to work with my project structure so I could upload my own apps, not just sample apps. Was able to successfully select my applet on the card and send a message, the applet checked its CLA and INS values, and returned them. Now have a better idea of how applet selection works.")
AddParaBeforeEnd("Wrote a small script to automate the process of compiling source code, converting to a .cap file, uploading it onto the card, and running it with a test host application.")
AddParaBeforeEnd("12/11 Wrote a test application that takes a byte string via an apdu reading “Hello World”, storing it, and returning the string upon a later apdu request.")

# Move the "_GoBack" bookmark from its original location (end of the
# "Created a new GitHub repository..." paragraph) to the end of the text
# of the new final content paragraph, then remove the leftover empty
# trailing paragraph so the new content paragraph becomes the last one.

$secondLast = $d.Paragraphs($d.Paragraphs.Count - 1)
$endOfText = $secondLast.Range.End - 1

# Directly placing a bookmark exactly at a paragraph-end offset is unreliable,
# so nudge it off that boundary with a temporary placeholder character,
# add the bookmark, then remove the placeholder again.
$insRng = $d.Range($endOfText, $endOfText)
$insRng.InsertAfter("Z")
$bmRng = $d.Range($endOfText, $endOfText)
$d.Bookmarks.Add("_GoBack", $bmRng)
$placeholderRng = $d.Range($endOfText, $endOfText + 1)
$placeholderRng.Delete()

# Remove the trailing empty paragraph left over from the original document.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$removeRng = $d.Range($lastPara.Range.Start - 1, $lastPara.Range.End)
$removeRng.Delete()
